$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: rows 2-11 (quarters Q0-Q9) had their
# underlying source-quarter shifted by one, so the recomputed
# ME/MAE/MSE/RMSE/SE statistics change, and N (col G) drops by 9
# uniformly (the sample-size impact of the reindex).

$ws.Range("B2").Value = 0.1091893978781611
$ws.Range("C2").Value = 2.452823248062764
$ws.Range("D2").Value = 19.2641810097944
$ws.Range("E2").Value = 4.389097972225546
$ws.Range("F2").Value = 4.448262992299666
$ws.Range("G2").Value = 37

$ws.Range("B3").Value = 0.5977652040503706
$ws.Range("C3").Value = 1.904146959540057
$ws.Range("D3").Value = 17.98508047146703
$ws.Range("E3").Value = 4.240882039324724
$ws.Range("F3").Value = 4.258099040526254
$ws.Range("G3").Value = 36

$ws.Range("B4").Value = 0.4593734525495876
$ws.Range("C4").Value = 2.116234692474444
$ws.Range("D4").Value = 19.0453471993615
$ws.Range("E4").Value = 4.364097524043373
$ws.Range("F4").Value = 4.403211772525009
$ws.Range("G4").Value = 35

$ws.Range("B5").Value = 0.6755236604910957
$ws.Range("C5").Value = 2.051568815471395
$ws.Range("D5").Value = 19.18090159300843
$ws.Range("E5").Value = 4.379600620263043
$ws.Range("F5").Value = 4.392263718217674
$ws.Range("G5").Value = 34

$ws.Range("B6").Value = 0.4638331138298176
$ws.Range("C6").Value = 2.158697177759521
$ws.Range("D6").Value = 20.1550692283716
$ws.Range("E6").Value = 4.489439745488473
$ws.Range("F6").Value = 4.534650022118768
$ws.Range("G6").Value = 33

$ws.Range("B7").Value = 0.7503312117146663
$ws.Range("C7").Value = 2.109841542719643
$ws.Range("D7").Value = 20.06410736424963
$ws.Range("E7").Value = 4.479297641846279
$ws.Range("F7").Value = 4.486666748889266
$ws.Range("G7").Value = 32

$ws.Range("B8").Value = 0.5544417764658882
$ws.Range("C8").Value = 2.186465380261041
$ws.Range("D8").Value = 21.15812832357595
$ws.Range("E8").Value = 4.599796552411417
$ws.Range("F8").Value = 4.641739623038775
$ws.Range("G8").Value = 31

$ws.Range("B9").Value = 0.7778321690867117
$ws.Range("C9").Value = 2.129509332268117
$ws.Range("D9").Value = 21.45336215612709
$ws.Range("E9").Value = 4.631777429467773
$ws.Range("F9").Value = 4.644055073279091
$ws.Range("G9").Value = 30

$ws.Range("B10").Value = 0.5115326828076273
$ws.Range("C10").Value = 2.258106810088389
$ws.Range("D10").Value = 22.44830329687462
$ws.Range("E10").Value = 4.737964045544733
$ws.Range("F10").Value = 4.793643450026643
$ws.Range("G10").Value = 29

$ws.Range("B11").Value = 0.8209835931991184
$ws.Range("C11").Value = 2.178446637574815
$ws.Range("D11").Value = 23.03067534287386
$ws.Range("E11").Value = 4.799028583252433
$ws.Range("F11").Value = 4.815047847583534
$ws.Range("G11").Value = 28
